$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 / Row 20: ShibaInu and Uniswap swap list positions, with updated price/volume ---
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0922"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
# D20 new value "6.27" parses as a number -- force text format so it is stored like the other price cells
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("E20").Value = "  -3.08%  "

# --- Remaining Price (D) / Volume(1h) (E) cell updates ---
$ws.Range("D2").Value = "40.914.32"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "2.415.86"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.83"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.60"
$ws.Range("E6").Value = "  -4.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.536"
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0832"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.85"
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "2.790.75"
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.69"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "2.419.89"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.773"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "40.885.71"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.01"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.93"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.29"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.68"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.02"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.21"
$ws.Range("E30").Value = "  -6.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.54"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0743"
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("E35").Value = "  -4.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.91"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.69"
$ws.Range("E37").Value = "  -4.16%  "
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -4.67%  "
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("D43").Value = "1.989.49"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.72"
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0273"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.40"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "2.658.63"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.13"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.21"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.52"
$ws.Range("E51").Value = "  -1.59%  "
